# Edit: add new base (F) functions `assertVarsNotPresent(vars)` / `assertVarsPresent(vars)`
# and new web (AA) function `assertLocation(search)` to the hidden '#system' sheet,
# keeping the existing alphabetically-sorted lists intact, then widen the
# corresponding defined names ("base", "web") to cover the extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---- Column F ("base" named range): full, alphabetically-sorted list ----
$baseValues = @(
    "appendText(var,appendWith)",
    "assertArrayContain(array,expected)",
    "assertArrayEqual(array1,array2,exactOrder)",
    "assertArrayNotContain(array,unexpected)",
    "assertContains(text,substring)",
    "assertCount(text,regex,expects)",
    "assertEmpty(text)",
    "assertEndsWith(text,suffix)",
    "assertEqual(expected,actual)",
    "assertMatch(text,regex)",
    "assertNotContain(text,substring)",
    "assertNotEmpty(text)",
    "assertNotEqual(expected,actual)",
    "assertStartsWith(text,prefix)",
    "assertTextOrder(var,descending)",
    "assertVarNotPresent(var)",
    "assertVarPresent(var)",
    "assertVarsNotPresent(vars)",
    "assertVarsPresent(vars)",
    "clear(vars)",
    "clearClipboard()",
    "copyFromClipboard(var)",
    "copyIntoClipboard(text)",
    "failImmediate(text)",
    "incrementChar(var,amount,config)",
    "macro(file,sheet,name)",
    "macroFlex(macro,input,output)",
    "outputToCloud(resource)",
    "prependText(var,prependWith)",
    "repeatUntil(steps,maxWaitMs)",
    "save(var,value)",
    "saveCount(text,regex,saveVar)",
    "saveMatches(text,regex,saveVar)",
    "saveReplace(text,regex,replace,saveVar)",
    "saveVariablesByPrefix(var,prefix)",
    "saveVariablesByRegex(var,regex)",
    "section(steps)",
    "split(text,delim,saveVar)",
    "startRecording()",
    "stopRecording()",
    "substringAfter(text,delim,saveVar)",
    "substringBefore(text,delim,saveVar)",
    "substringBetween(text,start,end,saveVar)",
    "verbose(text)",
    "waitFor(waitMs)",
    "waitForCondition(conditions,maxWaitMs)"
)

for ($i = 0; $i -lt $baseValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $baseValues[$i]
}

# ---- Column AA ("web" named range): full, alphabetically-sorted list ----
$webValues = @(
    "assertAndClick(locator,label)",
    "assertAttribute(locator,attrName,value)",
    "assertAttributeContain(locator,attrName,contains)",
    "assertAttributeNotContain(locator,attrName,contains)",
    "assertAttributeNotPresent(locator,attrName)",
    "assertAttributePresent(locator,attrName)",
    "assertChecked(locator)",
    "assertContainCount(locator,text,count)",
    "assertCssNotPresent(locator,property)",
    "assertCssPresent(locator,property,value)",
    "assertElementByAttributes(nameValues)",
    "assertElementByText(locator,text)",
    "assertElementCount(locator,count)",
    "assertElementDisabled(locator)",
    "assertElementEnabled(locator)",
    "assertElementNotPresent(locator)",
    "assertElementPresent(locator)",
    "assertElementsPresent(prefix)",
    "assertFocus(locator)",
    "assertFrameCount(count)",
    "assertFramePresent(frameName)",
    "assertIECompatMode()",
    "assertIENativeMode()",
    "assertLinkByLabel(label)",
    "assertLocation(search)",
    "assertMultiSelect(locator)",
    "assertNotChecked(locator)",
    "assertNotFocus(locator)",
    "assertNotText(locator,text)",
    "assertNotVisible(locator)",
    "assertOneMatch(locator)",
    "assertScrollbarHNotPresent(locator)",
    "assertScrollbarHPresent(locator)",
    "assertScrollbarVNotPresent(locator)",
    "assertScrollbarVPresent(locator)",
    "assertSingleSelect(locator)",
    "assertTable(locator,row,column,text)",
    "assertText(locator,text)",
    "assertTextContains(locator,text)",
    "assertTextCount(locator,text,count)",
    "assertTextList(locator,list,ignoreOrder)",
    "assertTextMatches(text,minMatch,scrollTo)",
    "assertTextNotContain(locator,text)",
    "assertTextNotPresent(text)",
    "assertTextOrder(locator,descending)",
    "assertTextPresent(text)",
    "assertTitle(text)",
    "assertValue(locator,value)",
    "assertValueOrder(locator,descending)",
    "assertVisible(locator)",
    "checkAll(locator,waitMs)",
    "clearLocalStorage()",
    "click(locator)",
    "clickAll(locator)",
    "clickAndWait(locator,waitMs)",
    "clickByLabel(label)",
    "clickByLabelAndWait(label,waitMs)",
    "clickOffset(locator,x,y)",
    "clickWithKeys(locator,keys)",
    "close()",
    "closeAll()",
    "deselect(locator,text)",
    "deselectMulti(locator,array)",
    "dismissInvalidCert()",
    "dismissInvalidCertPopup()",
    "doubleClick(locator)",
    "doubleClickAndWait(locator,waitMs)",
    "doubleClickByLabel(label)",
    "doubleClickByLabelAndWait(label,waitMs)",
    "dragAndDrop(fromLocator,toLocator)",
    "dragTo(fromLocator,xOffset,yOffset)",
    "editLocalStorage(key,value)",
    "executeScript(var,script)",
    "focus(locator)",
    "goBack()",
    "goBackAndWait()",
    "maximizeWindow()",
    "mouseOver(locator)",
    "moveTo(x,y)",
    "open(url)",
    "openAndWait(url,waitMs)",
    "openHttpBasic(url,username,password)",
    "openIgnoreTimeout(url)",
    "refresh()",
    "refreshAndWait()",
    "resizeWindow(width,height)",
    "rightClick(locator)",
    "saveAllWindowIds(var)",
    "saveAllWindowNames(var)",
    "saveAttribute(var,locator,attrName)",
    "saveAttributeList(var,locator,attrName)",
    "saveBrowserVersion(var)",
    "saveCount(var,locator)",
    "saveDivsAsCsv(headers,rows,cells,nextPage,file)",
    "saveElement(var,locator)",
    "saveElements(var,locator)",
    "saveInfiniteDivsAsCsv(config,file)",
    "saveInfiniteTableAsCsv(config,file)",
    "saveLocalStorage(var,key)",
    "saveLocation(var)",
    "savePageAs(var,sessionIdName,url)",
    "savePageAsFile(sessionIdName,url,file)",
    "saveSelectedText(var,locator)",
    "saveSelectedValue(var,locator)",
    "saveTableAsCsv(locator,nextPageLocator,file)",
    "saveText(var,locator)",
    "saveTextArray(var,locator)",
    "saveTextSubstringAfter(var,locator,delim)",
    "saveTextSubstringBefore(var,locator,delim)",
    "saveTextSubstringBetween(var,locator,start,end)",
    "saveTitle(var)",
    "saveValue(var,locator)",
    "saveValues(var,locator)",
    "screenshot(file,locator,removeFixed)",
    "screenshotInFull(file,timeout,removeFixed)",
    "scrollElement(locator,xOffset,yOffset)",
    "scrollLeft(locator,pixel)",
    "scrollPage(xOffset,yOffset)",
    "scrollRight(locator,pixel)",
    "scrollTo(locator)",
    "select(locator,text)",
    "selectAllOptions(locator)",
    "selectFrame(locator)",
    "selectMulti(locator,array)",
    "selectMultiByValue(locator,array)",
    "selectMultiOptions(locator)",
    "selectText(locator)",
    "selectWindow(winId)",
    "selectWindowAndWait(winId,waitMs)",
    "selectWindowByIndex(index)",
    "selectWindowByIndexAndWait(index,waitMs)",
    "switchBrowser(profile,config)",
    "toggleSelections(locator)",
    "type(locator,value)",
    "typeKeys(locator,value)",
    "uncheckAll(locator,waitMs)",
    "unselectAllText()",
    "updateAttribute(locator,attrName,value)",
    "upload(fieldLocator,file)",
    "verifyContainText(locator,text)",
    "verifyText(locator,text)",
    "wait(waitMs)",
    "waitForElementPresent(locator,waitMs)",
    "waitForElementsPresent(locators)",
    "waitForPopUp(winId,waitMs)",
    "waitForTextPresent(text)",
    "waitForTitle(text)",
    "waitUntilDisabled(locator,waitMs)",
    "waitUntilEnabled(locator,waitMs)",
    "waitUntilHidden(locator,waitMs)",
    "waitUntilVisible(locator,waitMs)"
)

for ($i = 0; $i -lt $webValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 27).Value = $webValues[$i]
}

# ---- Widen the defined names to cover the newly-added rows ----
$baseLastRow = 1 + $baseValues.Length
$webLastRow = 1 + $webValues.Length

$baseName = $wb.Names.Item("base")
$baseName.RefersTo = "='#system'!`$F`$2:`$F`$$($baseLastRow)"

$webName = $wb.Names.Item("web")
$webName.RefersTo = "='#system'!`$AA`$2:`$AA`$$($webLastRow)"

"base -> " + $baseName.RefersTo()
"web -> " + $webName.RefersTo()
